$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row at position 53. This pushes the old blank spacer
#    row (previously row 53, cells B53/C53 empty) down to become row 54,
#    and creates a fresh (duplicate-formatted) row 53 that we will fill in.
# ---------------------------------------------------------------------------
$ws.Rows("53:53").Insert()

# ---------------------------------------------------------------------------
# 2. Populate the new row 53 with the DELETE endpoint row that previously
#    lived in row 52 (text only - formatting was already copied down by the
#    row Insert above, matching the B52/C52 "s=1" style).
# ---------------------------------------------------------------------------
$ws.Range("B53").Value = "DELETE"
$ws.Range("C53").Value = "/users/{user-id}/accounts/{user-account-id}"

# ---------------------------------------------------------------------------
# 3. Change row 52 from DELETE to PUT (still same path).
# ---------------------------------------------------------------------------
$ws.Range("B52").Value = "PUT"
$ws.Range("C52").Value = "/users/{user-id}/accounts/{user-account-id}"

# ---------------------------------------------------------------------------
# 4. Mark column A for rows 2, 50, 51, 52, 53 with the same "highlight group"
#    formatting already used elsewhere in column A (copy format from A3,
#    which carries that style).
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("A50").PasteSpecial(-4122) | Out-Null
$ws.Range("A51").PasteSpecial(-4122) | Out-Null
$ws.Range("A52").PasteSpecial(-4122) | Out-Null
$ws.Range("A53").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5. Give B52 (now "PUT") its own distinct format: the same font/border
#    treatment as the rest of column B plus an explicit (blank) fill, which
#    is what the source workbook shows for this particular cell.
# ---------------------------------------------------------------------------
$b52 = $ws.Range("B52")
$b52.Borders.LineStyle = 1
$b52.Borders.Color = 16777215
$b52.Interior.Color = 16777215
$b52.Font.ThemeColor = 1
$b52.Font.TintAndShade = 0

# ---------------------------------------------------------------------------
# 6. Update the sheet view to match where the author left the selection.
# ---------------------------------------------------------------------------
$activeWindow = $excel.ActiveWindow
$activeWindow.ScrollRow = 22
$ws.Range("C56").Select()

Write-Host "Edit complete."
